# Auto-generated script applying the Midgardsormr_Profits market-data refresh
# (commit: "chore: update Sheets via scheduled runner")
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2285.9092
$ws.Range("I17").Value = 1300
$ws.Range("J17").Value = 2505
$ws.Range("K17").Value = 3900
$ws.Range("L17").Value = 7515
$ws.Range("M17").Value = -3732
$ws.Range("N17").Value = -7851
$ws.Range("H18").Value = 1639
$ws.Range("I18").Value = 1457.1428
$ws.Range("K18").Value = 1457.1428
$ws.Range("M18").Value = -1173.1428
$ws.Range("H33").Value = 859.95654
$ws.Range("J33").Value = 995.5
$ws.Range("L33").Value = 995.5
$ws.Range("N33").Value = -1453.5
$ws.Range("H64").Value = 8116.0586
$ws.Range("I64").Value = 4027
$ws.Range("K64").Value = 4027
$ws.Range("M64").Value = -3779
$ws.Range("H67").Value = 8116.0586
$ws.Range("I67").Value = 4027
$ws.Range("K67").Value = 4027
$ws.Range("M67").Value = -3169
$ws.Range("H76").Value = 4966.2856
$ws.Range("I76").Value = 4129.6665
$ws.Range("K76").Value = 4129.6665
$ws.Range("M76").Value = -3814.6665
$ws.Range("H79").Value = 4966.2856
$ws.Range("I79").Value = 4129.6665
$ws.Range("K79").Value = 4129.6665
$ws.Range("M79").Value = -3037.6665
$ws.Range("H98").Value = 2734.3
$ws.Range("I98").Value = 1435.2858
$ws.Range("J98").Value = 5765.3335
$ws.Range("K98").Value = 1435.2858
$ws.Range("L98").Value = 5765.3335
$ws.Range("M98").Value = 62.71419999999989
$ws.Range("N98").Value = -8761.333500000001
$ws.Range("H111").Value = 949.5
$ws.Range("I111").Value = 949.5
$ws.Range("K111").Value = 2848.5
$ws.Range("M111").Value = 218.5
$ws.Range("H122").Value = 2734.3
$ws.Range("I122").Value = 1435.2858
$ws.Range("J122").Value = 5765.3335
$ws.Range("K122").Value = 4305.857400000001
$ws.Range("L122").Value = 17296.0005
$ws.Range("M122").Value = -1855.857400000001
$ws.Range("N122").Value = -22196.0005
$ws.Range("H125").Value = 1646.75
$ws.Range("I125").Value = 799.5
$ws.Range("K125").Value = 7195.5
$ws.Range("M125").Value = -4735.5
$ws.Range("H131").Value = 3927.476
$ws.Range("I131").Value = 3656.8948
$ws.Range("K131").Value = 10970.6844
$ws.Range("M131").Value = -5930.6844
$ws.Range("H137").Value = 16228.869
$ws.Range("I137").Value = 63080.8
$ws.Range("J137").Value = 3214.4443
$ws.Range("K137").Value = 189242.4
$ws.Range("L137").Value = 9643.332900000001
$ws.Range("M137").Value = -186692.4
$ws.Range("N137").Value = -14743.3329
$ws.Range("H138").Value = 3720.9
$ws.Range("J138").Value = 4117
$ws.Range("L138").Value = 12351
$ws.Range("N138").Value = -22631
$ws.Range("H141").Value = 1157.2858
$ws.Range("I141").Value = 1157.2858
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 3471.8574
$ws.Range("L141").Value = 0
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = 1708.1426

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 75527.5
$ws.Range("I43").Value = 25000
$ws.Range("J43").Value = 92370
$ws.Range("K43").Value = 25000
$ws.Range("L43").Value = 92370
$ws.Range("M43").Value = -24687
$ws.Range("N43").Value = -92996
$ws.Range("H45").Value = 2930.08
$ws.Range("I45").Value = 2113.75
$ws.Range("J45").Value = 4381.3335
$ws.Range("K45").Value = 2113.75
$ws.Range("L45").Value = 4381.3335
$ws.Range("M45").Value = -1736.75
$ws.Range("N45").Value = -5135.3335
$ws.Range("H96").Value = 12500
$ws.Range("J96").Value = 12500
$ws.Range("L96").Value = 12500
$ws.Range("N96").Value = -17992
$ws.Range("H122").Value = 2272.75
$ws.Range("I122").Value = 2163.2144
$ws.Range("K122").Value = 6489.6432
$ws.Range("M122").Value = -4039.6432
$ws.Range("H132").Value = 1075.25
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5083.1113
$ws.Range("I134").Value = 5083.1113
$ws.Range("K134").Value = 15249.3339
$ws.Range("M134").Value = -12714.3339

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1722
$ws.Range("J16").Value = 1744.5
$ws.Range("L16").Value = 1744.5
$ws.Range("N16").Value = -2318.5
$ws.Range("H20").Value = 78390
$ws.Range("J20").Value = 78390
$ws.Range("L20").Value = 78390
$ws.Range("N20").Value = -78862
$ws.Range("H30").Value = 78390
$ws.Range("J30").Value = 78390
$ws.Range("L30").Value = 78390
$ws.Range("N30").Value = -78572
$ws.Range("H31").Value = 5557440.5
$ws.Range("I31").Value = 6251495.5
$ws.Range("J31").Value = 5000
$ws.Range("K31").Value = 6251495.5
$ws.Range("L31").Value = 5000
$ws.Range("M31").Value = -6251200.5
$ws.Range("N31").Value = -5590
$ws.Range("H34").Value = 5557440.5
$ws.Range("I34").Value = 6251495.5
$ws.Range("J34").Value = 5000
$ws.Range("K34").Value = 6251495.5
$ws.Range("L34").Value = 5000
$ws.Range("M34").Value = -6251293.5
$ws.Range("N34").Value = -5404
$ws.Range("H68").Value = 35000
$ws.Range("J68").Value = 35000
$ws.Range("L68").Value = 35000
$ws.Range("N68").Value = -36498
$ws.Range("H71").Value = 35000
$ws.Range("J71").Value = 35000
$ws.Range("L71").Value = 105000
$ws.Range("N71").Value = -112488
$ws.Range("H113").Value = 1722
$ws.Range("J113").Value = 1744.5
$ws.Range("L113").Value = 1744.5
$ws.Range("N113").Value = -6084.5
$ws.Range("H122").Value = 1140
$ws.Range("I122").Value = 945.7143
$ws.Range("K122").Value = 2837.1429
$ws.Range("M122").Value = -387.1428999999998
$ws.Range("H128").Value = 78390
$ws.Range("J128").Value = 78390
$ws.Range("L128").Value = 78390
$ws.Range("N128").Value = -88350
$ws.Range("H132").Value = 30965.7
$ws.Range("I132").Value = 32548.684
$ws.Range("K132").Value = 97646.052
$ws.Range("M132").Value = -95116.052
$ws.Range("H134").Value = 1747.3334
$ws.Range("I134").Value = 1474.3438
$ws.Range("J134").Value = 3931.25
$ws.Range("K134").Value = 4423.0314
$ws.Range("L134").Value = 11793.75
$ws.Range("M134").Value = -1888.0314
$ws.Range("N134").Value = -16863.75
$ws.Range("H141").Value = 175372.61
$ws.Range("J141").Value = 175372.61
$ws.Range("L141").Value = 175372.61
$ws.Range("N141").Value = -185732.61

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 143468.67
$ws.Range("J131").Value = 1845.3043
$ws.Range("L131").Value = 5535.9129
$ws.Range("N131").Value = -15615.9129
$ws.Range("H138").Value = 10166.4
$ws.Range("I138").Value = 2029.5
$ws.Range("K138").Value = 6088.5
$ws.Range("M138").Value = -948.5
$ws.Range("H140").Value = 2449.7144
$ws.Range("I140").Value = 2522.7693
$ws.Range("J140").Value = 1500
$ws.Range("K140").Value = 7568.3079
$ws.Range("L140").Value = 4500
$ws.Range("M140").Value = -2388.3079
$ws.Range("N140").Value = -14860
$ws.Range("H141").Value = 4923.615
$ws.Range("I141").Value = 5097.4
$ws.Range("J141").Value = 4344.3335
$ws.Range("K141").Value = 15292.2
$ws.Range("L141").Value = 13033.0005
$ws.Range("M141").Value = -10112.2
$ws.Range("N141").Value = -23393.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2587.4167
$ws.Range("I126").Value = 1832.8334
$ws.Range("K126").Value = 5498.5002
$ws.Range("M126").Value = -3028.5002
$ws.Range("H132").Value = 2690.9688
$ws.Range("I132").Value = 2327
$ws.Range("K132").Value = 6981
$ws.Range("M132").Value = -4451

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2180.5
$ws.Range("I7").Value = 2180.5
$ws.Range("K7").Value = 2180.5
$ws.Range("M7").Value = -2068.5
$ws.Range("H40").Value = 3402.1538
$ws.Range("I40").Value = 3053.25
$ws.Range("K40").Value = 3053.25
$ws.Range("M40").Value = -2917.25
$ws.Range("H126").Value = 2180.5
$ws.Range("I126").Value = 2180.5
$ws.Range("K126").Value = 6541.5
$ws.Range("M126").Value = -4071.5
$ws.Range("H132").Value = 2327.4614
$ws.Range("I132").Value = 1854.9166
$ws.Range("K132").Value = 5564.7498
$ws.Range("M132").Value = -3034.7498

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 835.63635
$ws.Range("I107").Value = 632.55554
$ws.Range("K107").Value = 1897.66662
$ws.Range("M107").Value = 22.33338000000003
$ws.Range("H132").Value = 1841.6604
$ws.Range("I132").Value = 1894.66
$ws.Range("K132").Value = 5683.98
$ws.Range("M132").Value = -3153.98
$ws.Range("H136").Value = 18068.432
$ws.Range("I136").Value = 20111.637
$ws.Range("K136").Value = 60334.91099999999
$ws.Range("M136").Value = -57784.91099999999
